# Add Change Language Button
# Insert a new "LANG" row right under the header row (row 1), pushing
# the existing language rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts rows 2.. down to 3..)
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "LANG"
$ws.Range("B2").Value = "English"
$ws.Range("C2").Value = "Tiếng Việt"

$ws.Range("C4").Select()
